$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B43").Value = 38
$ws.Range("B44").Value = 36
$ws.Range("B45").Value = 50
$ws.Range("B46").Value = 60
$ws.Range("B47").Value = 60
$ws.Range("B48").Value = 54
$ws.Range("B49").Value = 47
$ws.Range("B50").Value = 52
